$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (2023-09-11 -> 2023-09-12, i.e. 45180 -> 45181) for every data row.
$ws.Range("C2:C454").Value = 45181
